# Update LR-pair TPM data (Fgf2-Fgfr2): recompute the sending/target cluster
# matrix with the new TPM-derived values and add the previously-missing
# "Resolving-Mac" sending-cluster rows (13 -> 17 data rows total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7144740000000001
$ws.Range("N2").Value = 2.143422
$ws.Range("O2").Value = 0.138796410342318
$ws.Range("P2").Value = 0.138796410342318
$ws.Range("Q2").Value = 0.52310094831
$ws.Range("R2").Value = 4.70790853479
$ws.Range("S2").Value = 0.007097382806405967
$ws.Range("T2").Value = 0.007097382806405967
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.140873
$ws.Range("N3").Value = 12.422619
$ws.Range("O3").Value = 0.8044215857867821
$ws.Range("P3").Value = 0.8044215857867821
$ws.Range("Q3").Value = 3.031733265495
$ws.Range("R3").Value = 27.285599389455
$ws.Range("S3").Value = 0.04113426217568546
$ws.Range("T3").Value = 0.04113426217568546
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2847646666666667
$ws.Range("N4").Value = 0.8542940000000001
$ws.Range("O4").Value = 0.05531945672713084
$ws.Range("P4").Value = 0.05531945672713083
$ws.Range("Q4").Value = 0.2084899760922223
$ws.Range("R4").Value = 1.87640978483
$ws.Range("S4").Value = 0.002828771724474126
$ws.Range("T4").Value = 0.002828771724474126
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.007528666666666667
$ws.Range("N5").Value = 0.022586
$ws.Range("O5").Value = 0.00146254714376898
$ws.Range("P5").Value = 0.00146254714376898
$ws.Range("Q5").Value = 0.005512100752222222
$ws.Range("R5").Value = 0.04960890677
$ws.Range("S5").Value = 0.00007478764707345786
$ws.Range("T5").Value = 0.00007478764707345785
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7144740000000001
$ws.Range("N6").Value = 2.143422
$ws.Range("O6").Value = 0.138796410342318
$ws.Range("P6").Value = 0.138796410342318
$ws.Range("Q6").Value = 7.330388447844001
$ws.Range("R6").Value = 65.97349603059601
$ws.Range("S6").Value = 0.09945799773846511
$ws.Range("T6").Value = 0.09945799773846511
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.140873
$ws.Range("N7").Value = 12.422619
$ws.Range("O7").Value = 0.8044215857867821
$ws.Range("P7").Value = 0.8044215857867821
$ws.Range("Q7").Value = 42.484691679738
$ws.Range("R7").Value = 382.362225117642
$ws.Range("S7").Value = 0.5764281659924241
$ws.Range("T7").Value = 0.5764281659924241
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2847646666666667
$ws.Range("N8").Value = 0.8542940000000001
$ws.Range("O8").Value = 0.05531945672713084
$ws.Range("P8").Value = 0.05531945672713083
$ws.Range("Q8").Value = 2.921639727810223
$ws.Range("R8").Value = 26.294757550292
$ws.Range("S8").Value = 0.03964052376059606
$ws.Range("T8").Value = 0.03964052376059606
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.007528666666666667
$ws.Range("N9").Value = 0.022586
$ws.Range("O9").Value = 0.00146254714376898
$ws.Range("P9").Value = 0.00146254714376898
$ws.Range("Q9").Value = 0.07724291039422222
$ws.Range("R9").Value = 0.695186193548
$ws.Range("S9").Value = 0.001048024298024828
$ws.Range("T9").Value = 0.001048024298024828
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7144740000000001
$ws.Range("N10").Value = 2.143422
$ws.Range("O10").Value = 0.138796410342318
$ws.Range("P10").Value = 0.138796410342318
$ws.Range("Q10").Value = 1.834769470158
$ws.Range("R10").Value = 16.512925231422
$ws.Range("S10").Value = 0.02489397377941828
$ws.Range("T10").Value = 0.02489397377941828
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Fgf2"
$ws.Range("C11").Value = "Fgfr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.140873
$ws.Range("N11").Value = 12.422619
$ws.Range("O11").Value = 0.8044215857867821
$ws.Range("P11").Value = 0.8044215857867821
$ws.Range("Q11").Value = 10.633763244291
$ws.Range("R11").Value = 95.703869198619
$ws.Range("S11").Value = 0.144277865794838
$ws.Range("T11").Value = 0.144277865794838
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Fgf2"
$ws.Range("C12").Value = "Fgfr2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2847646666666667
$ws.Range("N12").Value = 0.8542940000000001
$ws.Range("O12").Value = 0.05531945672713084
$ws.Range("P12").Value = 0.05531945672713083
$ws.Range("Q12").Value = 0.7312757589215557
$ws.Range("R12").Value = 6.581481830294001
$ws.Range("S12").Value = 0.009921878396281444
$ws.Range("T12").Value = 0.009921878396281442
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Fgf2"
$ws.Range("C13").Value = "Fgfr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.007528666666666667
$ws.Range("N13").Value = 0.022586
$ws.Range("O13").Value = 0.00146254714376898
$ws.Range("P13").Value = 0.00146254714376898
$ws.Range("Q13").Value = 0.01933361850955556
$ws.Range("R13").Value = 0.174002566586
$ws.Range("S13").Value = 0.0002623166561610086
$ws.Range("T13").Value = 0.0002623166561610085
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Fgf2"
$ws.Range("C14").Value = "Fgfr2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.7144740000000001
$ws.Range("N14").Value = 2.143422
$ws.Range("O14").Value = 0.138796410342318
$ws.Range("P14").Value = 0.138796410342318
$ws.Range("Q14").Value = 0.541502702496
$ws.Range("R14").Value = 4.873524322464
$ws.Range("S14").Value = 0.007347056018028644
$ws.Range("T14").Value = 0.007347056018028644
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Fgf2"
$ws.Range("C15").Value = "Fgfr2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.140873
$ws.Range("N15").Value = 12.422619
$ws.Range("O15").Value = 0.8044215857867821
$ws.Range("P15").Value = 0.8044215857867821
$ws.Range("Q15").Value = 3.138384210192
$ws.Range("R15").Value = 28.245457891728
$ws.Range("S15").Value = 0.04258129182383449
$ws.Range("T15").Value = 0.04258129182383449
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Fgf2"
$ws.Range("C16").Value = "Fgfr2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2847646666666667
$ws.Range("N16").Value = 0.8542940000000001
$ws.Range("O16").Value = 0.05531945672713084
$ws.Range("P16").Value = 0.05531945672713083
$ws.Range("Q16").Value = 0.2158242799253333
$ws.Range("R16").Value = 1.942418519328
$ws.Range("S16").Value = 0.002928282845779209
$ws.Range("T16").Value = 0.002928282845779208
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Fgf2"
$ws.Range("C17").Value = "Fgfr2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.007528666666666667
$ws.Range("N17").Value = 0.022586
$ws.Range("O17").Value = 0.00146254714376898
$ws.Range("P17").Value = 0.00146254714376898
$ws.Range("Q17").Value = 0.005706006581333333
$ws.Range("R17").Value = 0.051354059232
$ws.Range("S17").Value = 0.00007741854250968542
$ws.Range("T17").Value = 0.00007741854250968541
